$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Re-format the header + existing data rows as Text (keeps fonts/fills, adds "@") ---
$ws.Range("A1:K1").NumberFormat = "@"
$ws.Range("A2:K2").NumberFormat = "@"
$ws.Range("A3:K3").NumberFormat = "@"
$ws.Range("B5").NumberFormat = "@"

# --- 2. Fix up A2 (was a numeric CCCD, must become the zero-padded text version) ---
$ws.Range("A2").Value = "001090001234"

# --- 3. New shareholder row (row 4) : Le Van D ---
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "001090001234"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "001090009012"
$ws.Range("C4").Value = 300
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "UQ-002"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2026-01-21"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "Uy quyen toan phan"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "Le Van D"
$ws.Range("I4").NumberFormat = "m/d/yy"
$ws.Range("I4").Value = "1/2/2022"

# --- 4. Newly-imported "ngay cap" for existing uy quyen rows 2-3 ---
$ws.Range("I2").Value = "22/12/2023"
$ws.Range("I3").Value = "22/12/2023"

# --- 5. Newly-imported "dia chi" column ---
$ws.Range("J2").Value = "Hà Nam"
$ws.Range("J3").Value = "Hà Nam"

# --- 6. Newly-imported "noi cap" column ---
$ws.Range("K2").Value = "Ninh Bình"
$ws.Range("K3").Value = "Thái Bình"

# --- 6b. Row 4 "noi cap" reuses the value just introduced above, default (General) format ---
$ws.Range("K4").Value = "Thái Bình"

# --- 7. Move the active selection the way the author left it ---
$ws.Range("G6").Select

$ws.PageSetup.Orientation = 1
